$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00001292064567892659
$ws.Range("C2").Value = 0.00006240767534437808
$ws.Range("D2").Value = 22.3905356188092
$ws.Range("E2").Value = 1133.036916526867
$ws.Range("G2").Value = 1155.427527473998
